# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple Leve-profit worksheets, refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 308.125
$ws.Range("I4").Value = 308.125
$ws.Range("K4").Value = 308.125
$ws.Range("M4").Value = -194.125

$ws.Range("H11").Value = 1002.7143
$ws.Range("I11").Value = 1002.7143
$ws.Range("K11").Value = 1002.7143
$ws.Range("M11").Value = -862.7143

$ws.Range("H12").Value = 28476.143
$ws.Range("I12").Value = 184
$ws.Range("J12").Value = 66199
$ws.Range("K12").Value = 184
$ws.Range("L12").Value = 66199
$ws.Range("M12").Value = -14
$ws.Range("N12").Value = -66539

$ws.Range("H33").Value = 434.57693
$ws.Range("I33").Value = 267.55554
$ws.Range("K33").Value = 267.55554
$ws.Range("M33").Value = -38.55554000000001

$ws.Range("H100").Value = 3277.4614
$ws.Range("I100").Value = 2039
$ws.Range("J100").Value = 4339
$ws.Range("K100").Value = 2039
$ws.Range("L100").Value = 4339
$ws.Range("M100").Value = -1498
$ws.Range("N100").Value = -5421

$ws.Range("H113").Value = 2255
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2255
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = 2255
$ws.Range("N113").Value = -8763
$ws.Range("L113").ClearContents()

$ws.Range("H116").Value = 11499.777
$ws.Range("I116").Value = 27499.334
$ws.Range("K116").Value = 27499.334
$ws.Range("M116").Value = -24057.334

$ws.Range("H121").Value = 4485.467
$ws.Range("J121").Value = 4485.467
$ws.Range("L121").Value = 13456.401
$ws.Range("N121").Value = -16950.401

$ws.Range("H132").Value = 7846.778
$ws.Range("I132").Value = 8816.134
$ws.Range("K132").Value = 26448.402
$ws.Range("M132").Value = -23918.402

$ws.Range("H141").Value = 3587.7144
$ws.Range("I141").Value = 3587.7144
$ws.Range("K141").Value = 10763.1432
$ws.Range("M141").Value = -5583.143199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3693.963
$ws.Range("I32").Value = 1793.9584
$ws.Range("K32").Value = 1793.9584
$ws.Range("M32").Value = -1506.9584

$ws.Range("H61").Value = 3181.5625
$ws.Range("I61").Value = 2850.4285
$ws.Range("K61").Value = 2850.4285
$ws.Range("M61").Value = -2638.4285

$ws.Range("H74").Value = 2363.375
$ws.Range("I74").Value = 2139.6155
$ws.Range("K74").Value = 2139.6155
$ws.Range("M74").Value = -1265.6155

$ws.Range("H77").Value = 2363.375
$ws.Range("I77").Value = 2139.6155
$ws.Range("K77").Value = 10698.0775
$ws.Range("M77").Value = -6330.077499999999

$ws.Range("H97").Value = 1215.6774
$ws.Range("I97").Value = 1063.4584
$ws.Range("K97").Value = 1063.4584
$ws.Range("M97").Value = -567.4584

$ws.Range("H123").Value = 100326.336
$ws.Range("J123").Value = 100326.336
$ws.Range("L123").Value = 100326.336
$ws.Range("N123").Value = -110126.336

$ws.Range("H132").Value = 2636.2856
$ws.Range("I132").Value = 2366.5642
$ws.Range("J132").Value = 3688.2
$ws.Range("K132").Value = 7099.692599999999
$ws.Range("L132").Value = 11064.6
$ws.Range("M132").Value = -4569.692599999999
$ws.Range("N132").Value = -16124.6

$ws.Range("H134").Value = 54999
$ws.Range("J134").Value = 54999
$ws.Range("L134").Value = 54999
$ws.Range("N134").Value = -65139

$ws.Range("H136").Value = 3181.5625
$ws.Range("I136").Value = 2850.4285
$ws.Range("K136").Value = 8551.2855
$ws.Range("M136").Value = -6001.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2200.5557
$ws.Range("I105").Value = 1885.8
$ws.Range("K105").Value = 1885.8
$ws.Range("M105").Value = -138.8

$ws.Range("H126").Value = 111998.5
$ws.Range("J126").Value = 111998.5
$ws.Range("L126").Value = 111998.5
$ws.Range("N126").Value = -121878.5

$ws.Range("H133").Value = 120001
$ws.Range("J133").Value = 120001
$ws.Range("L133").Value = 120001
$ws.Range("N133").Value = -130121

$ws.Range("H134").Value = 15875591
$ws.Range("I134").Value = 2237.0833
$ws.Range("J134").Value = 37040064
$ws.Range("K134").Value = 6711.249899999999
$ws.Range("L134").Value = 111120192
$ws.Range("M134").Value = -4176.249899999999
$ws.Range("N134").Value = -111125262

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2252.6667
$ws.Range("I5").Value = 379.5
$ws.Range("J5").Value = 5999
$ws.Range("K5").Value = 379.5
$ws.Range("L5").Value = 5999
$ws.Range("M5").Value = -267.5
$ws.Range("N5").Value = -6223

$ws.Range("H58").Value = 3068.6726
$ws.Range("I58").Value = 2752
$ws.Range("K58").Value = 2752
$ws.Range("M58").Value = -2549

$ws.Range("H105").Value = 2545.8
$ws.Range("I105").Value = 2265.6667
$ws.Range("K105").Value = 2265.6667
$ws.Range("M105").Value = -518.6667000000002

$ws.Range("H108").Value = 165000
$ws.Range("J108").Value = 165000
$ws.Range("L108").Value = 165000
$ws.Range("N108").Value = -172680

$ws.Range("H132").Value = 3049
$ws.Range("I132").Value = 3049
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9147
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -6617
$ws.Range("M132").ClearContents()

$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

$ws.Range("H136").Value = 3068.6726
$ws.Range("I136").Value = 2752
$ws.Range("K136").Value = 8256
$ws.Range("M136").Value = -5706

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1824.75
$ws.Range("I50").Value = 1159.8
$ws.Range("J50").Value = 2933
$ws.Range("K50").Value = 3479.4
$ws.Range("L50").Value = 8799
$ws.Range("M50").Value = -2998.4
$ws.Range("N50").Value = -9761

$ws.Range("H53").Value = 1824.75
$ws.Range("I53").Value = 1159.8
$ws.Range("J53").Value = 2933
$ws.Range("K53").Value = 3479.4
$ws.Range("L53").Value = 8799
$ws.Range("M53").Value = -2998.4
$ws.Range("N53").Value = -9761

$ws.Range("H129").Value = 872.2857
$ws.Range("I129").Value = 872.2857
$ws.Range("K129").Value = 2616.8571
$ws.Range("M129").Value = 2383.1429

$ws.Range("H132").Value = 266
$ws.Range("I132").Value = 266
$ws.Range("K132").Value = 2394
$ws.Range("M132").Value = 136

$ws.Range("H139").Value = 2888.8667
$ws.Range("I139").Value = 2562.0833
$ws.Range("K139").Value = 7686.249899999999
$ws.Range("M139").Value = -2546.249899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4589.6665
$ws.Range("I132").Value = 4877.5
$ws.Range("J132").Value = 4014
$ws.Range("K132").Value = 14632.5
$ws.Range("L132").Value = 12042
$ws.Range("M132").Value = -12102.5
$ws.Range("N132").Value = -17102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 31500
$ws.Range("J60").Value = 31500
$ws.Range("L60").Value = 31500
$ws.Range("N60").Value = -32518

$ws.Range("H93").Value = 43479200
$ws.Range("I93").Value = 66667330
$ws.Range("K93").Value = 66667330
$ws.Range("M93").Value = -66666082

$ws.Range("H122").Value = 23219.26
$ws.Range("I122").Value = 21165.5
$ws.Range("K122").Value = 63496.5
$ws.Range("M122").Value = -61046.5

$ws.Range("H125").Value = 92995.8
$ws.Range("J125").Value = 92995.8
$ws.Range("L125").Value = 92995.8
$ws.Range("N125").Value = -102835.8

$ws.Range("H132").Value = 2828.5
$ws.Range("I132").Value = 2495.8
$ws.Range("K132").Value = 7487.400000000001
$ws.Range("M132").Value = -4957.400000000001

$ws.Range("H136").Value = 4667.3335
$ws.Range("I136").Value = 4126
$ws.Range("K136").Value = 12378
$ws.Range("M136").Value = -9828

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1959.25
$ws.Range("I100").Value = 2557
$ws.Range("J100").Value = 963
$ws.Range("K100").Value = 5114
$ws.Range("L100").Value = 1926
$ws.Range("M100").Value = -4573
$ws.Range("N100").Value = -3008

$ws.Range("H126").Value = 16596.6
$ws.Range("I126").Value = 16596.6
$ws.Range("K126").Value = 49789.8
$ws.Range("M126").Value = -47319.8

$ws.Range("H132").Value = 2386.2922
$ws.Range("I132").Value = 2329.7192
$ws.Range("K132").Value = 6989.1576
$ws.Range("M132").Value = -4459.1576
